# Added PeopleSoft IDs for TTL and collaborators
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# field_wbddh_dsttl_upi value
$ws.Range("B5").Value = 15062

# field_wbddh_collaborator_upi value
$ws.Range("B11").Value = "23715, 54524"

$ws.Range("B11").Select()
